$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.983.39"
$ws.Range("E2").Value = "  +6.55%  "
$ws.Range("D3").Value = "3.013.16"
$ws.Range("E3").Value = "  +3.91%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "584.15"
$ws.Range("E5").Value = "  +2.67%  "
$ws.Range("D6").Value = "162.82"
$ws.Range("E6").Value = "  +13.08%  "
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("D8").Value = "3.008.38"
$ws.Range("E8").Value = "  +3.83%  "
$ws.Range("E9").Value = "  +3.21%  "
$ws.Range("D10").Value = "6.94"
$ws.Range("E10").Value = "  +0.52%  "
$ws.Range("D11").Value = "0.157"
$ws.Range("E11").Value = "  +7.31%  "
$ws.Range("D12").Value = "0.459"
$ws.Range("E12").Value = "  +6.49%  "
$ws.Range("D13").Value = "0.0000251"
$ws.Range("E13").Value = "  +8.79%  "
$ws.Range("D14").Value = "34.90"
$ws.Range("E14").Value = "  +8.28%  "
$ws.Range("D15").Value = "0.124"
$ws.Range("E15").Value = "  -0.86%  "
$ws.Range("D16").Value = "65.996.57"
$ws.Range("E16").Value = "  +6.65%  "
$ws.Range("D17").Value = "3.514.50"
$ws.Range("E17").Value = "  +3.96%  "
$ws.Range("D18").Value = "6.98"
$ws.Range("E18").Value = "  +7.25%  "
$ws.Range("D19").Value = "3.015.19"
$ws.Range("E19").Value = "  +3.74%  "
$ws.Range("D20").Value = "458.90"
$ws.Range("E20").Value = "  +6.31%  "
$ws.Range("D21").Value = "13.96"
$ws.Range("E21").Value = "  +8.02%  "
$ws.Range("D22").Value = "0.687"
$ws.Range("E22").Value = "  +5.26%  "
$ws.Range("E23").Value = "  +7.91%  "
$ws.Range("D24").Value = "82.57"
$ws.Range("E24").Value = "  +4.68%  "
$ws.Range("D25").Value = "2.33"
$ws.Range("E25").Value = "  +14.84%  "
$ws.Range("E26").Value = "  +3.70%  "
$ws.Range("D27").Value = "10.58"
$ws.Range("E27").Value = "  +4.81%  "
$ws.Range("E29").Value = "  +16.80%  "
$ws.Range("E30").Value = "  +15.78%  "
$ws.Range("E31").Value = "  +4.28%  "
$ws.Range("E32").Value = "  -6.63%  "
$ws.Range("D33").Value = "27.03"
$ws.Range("E33").Value = "  +5.54%  "
$ws.Range("E34").Value = "  +3.38%  "
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  -0.07%  "
$ws.Range("E36").Value = "  +4.20%  "
$ws.Range("D37").Value = "5.80"
$ws.Range("E37").Value = "  +7.74%  "
$ws.Range("D38").Value = "2.15"
$ws.Range("E38").Value = "  +12.19%  "
$ws.Range("D39").Value = "3.03"
$ws.Range("E39").Value = "  +6.36%  "
$ws.Range("D40").Value = "49.88"
$ws.Range("E40").Value = "  +2.22%  "
$ws.Range("D41").Value = "0.312"
$ws.Range("E41").Value = "  +16.04%  "
$ws.Range("E42").Value = "  +6.09%  "
$ws.Range("D43").Value = "43.85"
$ws.Range("E43").Value = "  +8.78%  "
$ws.Range("E44").Value = "  +4.36%  "
$ws.Range("D45").Value = "388.09"
$ws.Range("E45").Value = "  +12.03%  "
$ws.Range("E46").Value = "  +6.12%  "
$ws.Range("D47").Value = "2.795.63"
$ws.Range("E47").Value = "  +3.49%  "
$ws.Range("D48").Value = "135.35"
$ws.Range("E48").Value = "  +2.77%  "
$ws.Range("D50").Value = "24.00"
$ws.Range("E50").Value = "  +11.22%  "
$ws.Range("E51").Value = "  +4.14%  "
